# Insert a new weekly price record for Camote (Vega Central Mapocho de Santiago)
# at row 17, pushing the existing rows 17-67 down to 18-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..67 down to 18..68 by inserting a new row at 17.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new observation.
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44592
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100114002
$ws.Cells.Item(17, 7).Value = "Camote"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 1000
$ws.Cells.Item(17, 11).Value = 11000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 11600
$ws.Cells.Item(17, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 644
$ws.Cells.Item(17, 17).Value = 18
$ws.Cells.Item(17, 18).Value = "Hortaliza"
